$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "thuat toan" text marker cell (";") next to the existing data table.
$ws.Range("K11").Value = ";"

# Switch the workbook's default/Normal font from Arial to Calibri.
$wb.Styles.Item("Normal").Font.Name = "Calibri"

# Move the active selection to G5, as left by the author after the edit.
[void]$ws.Range("G5").Select()
